# Scheduled market-data refresh: update per-item average prices and the
# resulting leve price/profit figures across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets (columns H-N: currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2468.3333
$ws.Range("I62").Value = 2468.3333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2468.3333
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -1844.3333
$ws.Range("H65").Value = 2468.3333
$ws.Range("I65").Value = 2468.3333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 12341.6665
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -9221.666499999999
$ws.Range("H69").Value = 8392.5
$ws.Range("J69").Value = 9790
$ws.Range("L69").Value = 29370
$ws.Range("N69").Value = -31118
$ws.Range("H72").Value = 8392.5
$ws.Range("J72").Value = 9790
$ws.Range("L72").Value = 88110
$ws.Range("N72").Value = -96846
$ws.Range("H98").Value = 790
$ws.Range("I98").Value = 733.3333
$ws.Range("J98").Value = 960
$ws.Range("K98").Value = 733.3333
$ws.Range("L98").Value = 960
$ws.Range("M98").Value = 764.6667
$ws.Range("N98").Value = -3956
$ws.Range("H122").Value = 790
$ws.Range("I122").Value = 733.3333
$ws.Range("J122").Value = 960
$ws.Range("K122").Value = 2199.9999
$ws.Range("L122").Value = 2880
$ws.Range("M122").Value = 250.0001000000002
$ws.Range("N122").Value = -7780
$ws.Range("H137").Value = 1521.0364
$ws.Range("I137").Value = 1377.125
$ws.Range("J137").Value = 1632.4517
$ws.Range("K137").Value = 4131.375
$ws.Range("L137").Value = 4897.355100000001
$ws.Range("M137").Value = -1581.375
$ws.Range("N137").Value = -9997.355100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2299.125
$ws.Range("I88").Value = 1850
$ws.Range("J88").Value = 2448.8333
$ws.Range("K88").Value = 1850
$ws.Range("L88").Value = 2448.8333
$ws.Range("M88").Value = -1444
$ws.Range("N88").Value = -3260.8333
$ws.Range("H91").Value = 2299.125
$ws.Range("I91").Value = 1850
$ws.Range("J91").Value = 2448.8333
$ws.Range("K91").Value = 1850
$ws.Range("L91").Value = 2448.8333
$ws.Range("M91").Value = -446
$ws.Range("N91").Value = -5256.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 41022.785
$ws.Range("I86").Value = 123604
$ws.Range("J86").Value = 1905.3684
$ws.Range("K86").Value = 123604
$ws.Range("L86").Value = 1905.3684
$ws.Range("M86").Value = -122481
$ws.Range("N86").Value = -4151.3684
$ws.Range("H89").Value = 41022.785
$ws.Range("I89").Value = 123604
$ws.Range("J89").Value = 1905.3684
$ws.Range("K89").Value = 618020
$ws.Range("L89").Value = 9526.842000000001
$ws.Range("M89").Value = -612404
$ws.Range("N89").Value = -20758.842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18077.572
$ws.Range("I31").Value = 39620.77
$ws.Range("J31").Value = 2939.1082
$ws.Range("K31").Value = 39620.77
$ws.Range("L31").Value = 2939.1082
$ws.Range("M31").Value = -39325.77
$ws.Range("N31").Value = -3529.1082
$ws.Range("H34").Value = 18077.572
$ws.Range("I34").Value = 39620.77
$ws.Range("J34").Value = 2939.1082
$ws.Range("K34").Value = 39620.77
$ws.Range("L34").Value = 2939.1082
$ws.Range("M34").Value = -39418.77
$ws.Range("N34").Value = -3343.1082
$ws.Range("H94").Value = 1269.1538
$ws.Range("J94").Value = 1368.6
$ws.Range("L94").Value = 1368.6
$ws.Range("N94").Value = -2270.6
$ws.Range("H99").Value = 11523.786
$ws.Range("J99").Value = 14872.556
$ws.Range("L99").Value = 14872.556
$ws.Range("N99").Value = -17868.556
$ws.Range("H122").Value = 2708.2693
$ws.Range("I122").Value = 2557.2104
$ws.Range("J122").Value = 3118.2856
$ws.Range("K122").Value = 7671.6312
$ws.Range("L122").Value = 9354.856800000001
$ws.Range("M122").Value = -5221.6312
$ws.Range("N122").Value = -14254.8568
$ws.Range("H126").Value = 11523.786
$ws.Range("J126").Value = 14872.556
$ws.Range("L126").Value = 44617.66800000001
$ws.Range("N126").Value = -49557.66800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 200
$ws.Range("I6").Value = 200
$ws.Range("K6").Value = 600
$ws.Range("M6").Value = -487
$ws.Range("H107").Value = 268337.06
$ws.Range("I107").Value = 649.96155
$ws.Range("J107").Value = 536024.1
$ws.Range("K107").Value = 1949.88465
$ws.Range("L107").Value = 1608072.3
$ws.Range("M107").Value = -29.88464999999997
$ws.Range("N107").Value = -1611912.3
$ws.Range("H113").Value = 544.6957
$ws.Range("J113").Value = 553
$ws.Range("L113").Value = 1659
$ws.Range("N113").Value = -5999
$ws.Range("H132").Value = 2574.7307
$ws.Range("I132").Value = 2417.1667
$ws.Range("J132").Value = 2622
$ws.Range("K132").Value = 21754.5003
$ws.Range("L132").Value = 23598
$ws.Range("M132").Value = -19224.5003
$ws.Range("N132").Value = -28658

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 3861.5
$ws.Range("I99").Value = 2416
$ws.Range("J99").Value = 13980
$ws.Range("K99").Value = 2416
$ws.Range("L99").Value = 13980
$ws.Range("M99").Value = -170
$ws.Range("N99").Value = -18472
$ws.Range("H100").Value = 35352.5
$ws.Range("J100").Value = 35352.5
$ws.Range("L100").Value = 35352.5
$ws.Range("N100").Value = -37516.5
$ws.Range("H105").Value = 38990
$ws.Range("J105").Value = 38990
$ws.Range("L105").Value = 38990
$ws.Range("N105").Value = -45978

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1840.5217
$ws.Range("I7").Value = 1395.6364
$ws.Range("K7").Value = 1395.6364
$ws.Range("M7").Value = -1283.6364
$ws.Range("H61").Value = 2174.625
$ws.Range("I61").Value = 2082.8333
$ws.Range("J61").Value = 2450
$ws.Range("K61").Value = 2082.8333
$ws.Range("L61").Value = 2450
$ws.Range("M61").Value = -1880.8333
$ws.Range("N61").Value = -2854
$ws.Range("H113").Value = 2174.625
$ws.Range("I113").Value = 2082.8333
$ws.Range("J113").Value = 2450
$ws.Range("K113").Value = 2082.8333
$ws.Range("L113").Value = 2450
$ws.Range("M113").Value = 87.16670000000022
$ws.Range("N113").Value = -6790
$ws.Range("H126").Value = 1840.5217
$ws.Range("I126").Value = 1395.6364
$ws.Range("K126").Value = 4186.9092
$ws.Range("M126").Value = -1716.9092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3398.1936
$ws.Range("I132").Value = 4705.6875
$ws.Range("J132").Value = 2003.5333
$ws.Range("K132").Value = 14117.0625
$ws.Range("L132").Value = 6010.5999
$ws.Range("M132").Value = -11587.0625
$ws.Range("N132").Value = -11070.5999
